$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 560328.9399999999
$ws.Range("I64").Value = 772804.7
$ws.Range("J64").Value = 7892
$ws.Range("K64").Value = 772804.7
$ws.Range("L64").Value = 7892
$ws.Range("M64").Value = -772556.7
$ws.Range("N64").Value = -8388
$ws.Range("H67").Value = 560328.9399999999
$ws.Range("I67").Value = 772804.7
$ws.Range("J67").Value = 7892
$ws.Range("K67").Value = 772804.7
$ws.Range("L67").Value = 7892
$ws.Range("M67").Value = -771946.7
$ws.Range("N67").Value = -9608
$ws.Range("H86").Value = 1888.1154
$ws.Range("I86").Value = 1491.4546
$ws.Range("J86").Value = 2179
$ws.Range("K86").Value = 1491.4546
$ws.Range("L86").Value = 2179
$ws.Range("M86").Value = -368.4546
$ws.Range("N86").Value = -4425
$ws.Range("H89").Value = 1888.1154
$ws.Range("I89").Value = 1491.4546
$ws.Range("J89").Value = 2179
$ws.Range("K89").Value = 7457.273
$ws.Range("L89").Value = 10895
$ws.Range("M89").Value = -1841.273
$ws.Range("N89").Value = -22127
$ws.Range("H133").Value = 20000
$ws.Range("J133").Value = 20000
$ws.Range("L133").Value = 20000
$ws.Range("N133").Value = -30120
$ws.Range("H137").Value = 27028336
$ws.Range("I137").Value = 33334224
$ws.Range("J137").Value = 3100.2856
$ws.Range("K137").Value = 100002672
$ws.Range("L137").Value = 9300.856800000001
$ws.Range("M137").Value = -100000122
$ws.Range("N137").Value = -14400.8568

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1322.6428
$ws.Range("I45").Value = 1000.3333
$ws.Range("J45").Value = 1902.8
$ws.Range("K45").Value = 1000.3333
$ws.Range("L45").Value = 1902.8
$ws.Range("M45").Value = -623.3333
$ws.Range("N45").Value = -2656.8
$ws.Range("H61").Value = 3618.5
$ws.Range("I61").Value = 2196.923
$ws.Range("J61").Value = 4705.5884
$ws.Range("K61").Value = 2196.923
$ws.Range("L61").Value = 4705.5884
$ws.Range("M61").Value = -1984.923
$ws.Range("N61").Value = -5129.5884
$ws.Range("H63").Value = 43333.332
$ws.Range("I63").Value = 43333.332
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 43333.332
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -42647.332
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 43333.332
$ws.Range("I66").Value = 43333.332
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 216666.66
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -213234.66
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 4695.25
$ws.Range("I74").Value = 1241.25
$ws.Range("J74").Value = 11603.25
$ws.Range("K74").Value = 1241.25
$ws.Range("L74").Value = 11603.25
$ws.Range("M74").Value = -367.25
$ws.Range("N74").Value = -13351.25
$ws.Range("H77").Value = 4695.25
$ws.Range("I77").Value = 1241.25
$ws.Range("J77").Value = 11603.25
$ws.Range("K77").Value = 6206.25
$ws.Range("L77").Value = 58016.25
$ws.Range("M77").Value = -1838.25
$ws.Range("N77").Value = -66752.25
$ws.Range("H132").Value = 4463.9414
$ws.Range("I132").Value = 4714
$ws.Range("J132").Value = 4241.6665
$ws.Range("K132").Value = 14142
$ws.Range("L132").Value = 12724.9995
$ws.Range("M132").Value = -11612
$ws.Range("N132").Value = -17784.9995
$ws.Range("H136").Value = 3618.5
$ws.Range("I136").Value = 2196.923
$ws.Range("J136").Value = 4705.5884
$ws.Range("K136").Value = 6590.768999999999
$ws.Range("L136").Value = 14116.7652
$ws.Range("M136").Value = -4040.768999999999
$ws.Range("N136").Value = -19216.7652
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 34735.668
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 51103.5
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 51103.5
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -53349.5
$ws.Range("H89").Value = 34735.668
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 51103.5
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 255517.5
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -266749.5
$ws.Range("H134").Value = 2593.3684
$ws.Range("I134").Value = 1900.1818
$ws.Range("K134").Value = 5700.5454
$ws.Range("M134").Value = -3165.5454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1451.0952
$ws.Range("I31").Value = 1182.7894
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1182.7894
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -887.7893999999999
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 1451.0952
$ws.Range("I34").Value = 1182.7894
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1182.7894
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -980.7893999999999
$ws.Range("N34").Value = -4404
$ws.Range("H58").Value = 1798.697
$ws.Range("I58").Value = 1242.75
$ws.Range("J58").Value = 3281.2222
$ws.Range("K58").Value = 1242.75
$ws.Range("L58").Value = 3281.2222
$ws.Range("M58").Value = -1039.75
$ws.Range("N58").Value = -3687.2222
$ws.Range("H132").Value = 2525.0386
$ws.Range("I132").Value = 2035.7727
$ws.Range("J132").Value = 5216
$ws.Range("K132").Value = 6107.3181
$ws.Range("L132").Value = 15648
$ws.Range("M132").Value = -3577.3181
$ws.Range("N132").Value = -20708
$ws.Range("H134").Value = 2119.195
$ws.Range("I134").Value = 1079.5
$ws.Range("K134").Value = 3238.5
$ws.Range("M134").Value = -703.5
$ws.Range("H136").Value = 1798.697
$ws.Range("I136").Value = 1242.75
$ws.Range("J136").Value = 3281.2222
$ws.Range("K136").Value = 3728.25
$ws.Range("L136").Value = 9843.6666
$ws.Range("M136").Value = -1178.25
$ws.Range("N136").Value = -14943.6666

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 685.1579
$ws.Range("I23").Value = 1491.4286
$ws.Range("J23").Value = 214.83333
$ws.Range("K23").Value = 4474.2858
$ws.Range("L23").Value = 644.49999
$ws.Range("M23").Value = -4239.2858
$ws.Range("N23").Value = -1114.49999
$ws.Range("H97").Value = 777.4666999999999
$ws.Range("I97").Value = 889.5
$ws.Range("J97").Value = 649.4286
$ws.Range("K97").Value = 2668.5
$ws.Range("L97").Value = 1948.2858
$ws.Range("M97").Value = -2172.5
$ws.Range("N97").Value = -2940.2858
$ws.Range("H98").Value = 233.92857
$ws.Range("J98").Value = 233.25
$ws.Range("L98").Value = 699.75
$ws.Range("N98").Value = -3695.75
$ws.Range("H113").Value = 11364381
$ws.Range("I113").Value = 539.6667
$ws.Range("J113").Value = 13158672
$ws.Range("K113").Value = 1619.0001
$ws.Range("L113").Value = 39476016
$ws.Range("M113").Value = 550.9999
$ws.Range("N113").Value = -39480356
$ws.Range("H117").Value = 900
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 900
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 2700
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -9584

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2497.6875
$ws.Range("I132").Value = 2113.889
$ws.Range("J132").Value = 3649.0833
$ws.Range("K132").Value = 6341.667
$ws.Range("L132").Value = 10947.2499
$ws.Range("M132").Value = -3811.667
$ws.Range("N132").Value = -16007.2499

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2970.8235
$ws.Range("I7").Value = 1900.8
$ws.Range("J7").Value = 3416.6667
$ws.Range("K7").Value = 1900.8
$ws.Range("L7").Value = 3416.6667
$ws.Range("M7").Value = -1788.8
$ws.Range("N7").Value = -3640.6667
$ws.Range("H22").Value = 5908.8
$ws.Range("I22").Value = 511.75
$ws.Range("J22").Value = 9506.833000000001
$ws.Range("K22").Value = 511.75
$ws.Range("L22").Value = 9506.833000000001
$ws.Range("M22").Value = -216.75
$ws.Range("N22").Value = -10096.833
$ws.Range("H27").Value = 5908.8
$ws.Range("I27").Value = 511.75
$ws.Range("J27").Value = 9506.833000000001
$ws.Range("K27").Value = 511.75
$ws.Range("L27").Value = 9506.833000000001
$ws.Range("M27").Value = -404.75
$ws.Range("N27").Value = -9720.833000000001
$ws.Range("H122").Value = 3600
$ws.Range("I122").Value = 2800
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8400
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -5950
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 2970.8235
$ws.Range("I126").Value = 1900.8
$ws.Range("J126").Value = 3416.6667
$ws.Range("K126").Value = 5702.4
$ws.Range("L126").Value = 10250.0001
$ws.Range("M126").Value = -3232.4
$ws.Range("N126").Value = -15190.0001
$ws.Range("H132").Value = 3666.2036
$ws.Range("I132").Value = 2917.9722
$ws.Range("J132").Value = 5162.6665
$ws.Range("K132").Value = 8753.9166
$ws.Range("L132").Value = 15487.9995
$ws.Range("M132").Value = -6223.9166
$ws.Range("N132").Value = -20547.9995
$ws.Range("H136").Value = 3637.1592
$ws.Range("I136").Value = 2121.3438
$ws.Range("J136").Value = 7679.3335
$ws.Range("K136").Value = 6364.0314
$ws.Range("L136").Value = 23038.0005
$ws.Range("M136").Value = -3814.0314
$ws.Range("N136").Value = -28138.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 61576.293
$ws.Range("I122").Value = 127188.375
$ws.Range("J122").Value = 3254.4443
$ws.Range("K122").Value = 381565.125
$ws.Range("L122").Value = 9763.332900000001
$ws.Range("M122").Value = -379115.125
$ws.Range("N122").Value = -14663.3329
$ws.Range("H123").Value = 25561.285
$ws.Range("J123").Value = 25561.285
$ws.Range("L123").Value = 25561.285
$ws.Range("N123").Value = -35361.285
$ws.Range("H126").Value = 112667.22
$ws.Range("I126").Value = 168200.83
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 504602.49
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -502132.49
$ws.Range("N126").Value = -9740
$ws.Range("H132").Value = 17860342
$ws.Range("I132").Value = 22731016
$ws.Range("K132").Value = 68193048
$ws.Range("M132").Value = -68190518
$ws.Range("H136").Value = 9037295
$ws.Range("I136").Value = 13931500
$ws.Range("J136").Value = 1840.6154
$ws.Range("K136").Value = 41794500
$ws.Range("L136").Value = 5521.8462
$ws.Range("M136").Value = -41791950
$ws.Range("N136").Value = -10621.8462
